# Update the "Förändrad" (Changed) date column (C) for data rows 2-9
# from serial date 45224 (2023-10-25) to 45233 (2023-11-03).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($row = 2; $row -le 9; $row++) {
    $ws.Cells.Item($row, 3).Value = 45233
}
